$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.209.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5231"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06314"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.531"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.621.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8066"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.201.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.719"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.020"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1206"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.229"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05603"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.47%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.485"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.376"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.797"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.405"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5740"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01600"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.567"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8466"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.041.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.795.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₈106"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05315"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.046"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
